$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5501
$ws1.Range("F4").Value = 12136
$ws1.Range("F8").Value = 331
$ws1.Range("F9").Value = 1106

# Sheet "全部类型" (all types) - same events repeated, update same column
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5501
$ws4.Range("F6").Value = 12136
$ws4.Range("F12").Value = 331
$ws4.Range("F13").Value = 1106
